$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old "total" row (row 11) so its content can be rebuilt at row 27 ---
$ws.Range("A11:B11").ClearContents() | Out-Null

# --- Header row: add new "type" column ---
$ws.Range("E1").Value = "type"

# --- Row 2-9: existing rows, update "Assigned To" to Vaughan and add "type" = CSS ---
$ws.Range("C2").Value = "Vaughan"
$ws.Range("E2").Value = "CSS"

$ws.Range("C3").Value = "Vaughan"
$ws.Range("E3").Value = "CSS"

$ws.Range("C4").Value = "Vaughan"
$ws.Range("E4").Value = "CSS"

$ws.Range("C5").Value = "Vaughan"
$ws.Range("E5").Value = "CSS"

$ws.Range("C6").Value = "Vaughan"
$ws.Range("E6").Value = "CSS"

$ws.Range("C7").Value = "Vaughan"
$ws.Range("E7").Value = "CSS"

$ws.Range("C8").Value = "Vaughan"
$ws.Range("E8").Value = "CSS"

$ws.Range("C9").Value = "Vaughan"
$ws.Range("E9").Value = "CSS"

# --- Rows 10-17: HTML block (repeat of tasks) ---
$ws.Range("A10").Value = "title and searchbar"
$ws.Range("B10").Value = 0.5
$ws.Range("B10").HorizontalAlignment = -4152
$ws.Range("C10").Value = "Caelan"
$ws.Range("E10").Value = "HTML"

$ws.Range("A11").Value = "Profile"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Rowan"
$ws.Range("E11").Value = "HTML"

$ws.Range("A12").Value = "Competitive Overview"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "Caelan"
$ws.Range("E12").Value = "HTML"

$ws.Range("A13").Value = "graph integration"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "Rowan"
$ws.Range("E13").Value = "HTML"

$ws.Range("A14").Value = "champ mastery"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "Caelan"
$ws.Range("E14").Value = "HTML"

$ws.Range("A15").Value = "match history"
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "Rowan"
$ws.Range("E15").Value = "HTML"

$ws.Range("A16").Value = "bugchecking/fixing"
$ws.Range("B16").Value = 0.5
$ws.Range("C16").Value = "Caelan"
$ws.Range("E16").Value = "HTML"

$ws.Range("A17").Value = "finishing touches + testing + extra features"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Rowan"
$ws.Range("E17").Value = "HTML"

# --- Rows 18-25: JavaScript block (repeat of tasks) ---
$ws.Range("A18").Value = "title and searchbar"
$ws.Range("B18").Value = 0.5
$ws.Range("B18").HorizontalAlignment = -4152
$ws.Range("C18").Value = "Caelan"
$ws.Range("E18").Value = "JavaScript"

$ws.Range("A19").Value = "Profile"
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "Rowan"
$ws.Range("E19").Value = "JavaScript"

$ws.Range("A20").Value = "Competitive Overview"
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = "Caelan"
$ws.Range("E20").Value = "JavaScript"

$ws.Range("A21").Value = "graph integration"
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = "Rowan"
$ws.Range("E21").Value = "JavaScript"

$ws.Range("A22").Value = "champ mastery"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = "Caelan"
$ws.Range("E22").Value = "JavaScript"

$ws.Range("A23").Value = "match history"
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = "Rowan"
$ws.Range("E23").Value = "JavaScript"

$ws.Range("A24").Value = "bugchecking/fixing"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = "Caelan"
$ws.Range("E24").Value = "JavaScript"

$ws.Range("A25").Value = "finishing touches + testing + extra features"
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = "Rowan"
$ws.Range("E25").Value = "JavaScript"

# --- Row 27: new total row (row 26 intentionally left blank) ---
$ws.Range("A27").Value = "total"
$ws.Range("B27").Formula = "=SUM(B2:B26)"

# --- Update selection to match the final state ---
$ws.Range("B25").Select() | Out-Null
